$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three changed cell values
$ws.Range("C2").Value = "Donald Trump1"
$ws.Range("B3").Value = "Hannah Flores"
$ws.Range("K4").Value = " Trump1"

# Move the active selection to K4 (matches the saved view state)
$ws.Range("K4").Select()
